$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string text: remove the trailing period.
$ws.Range("B2").Value = "Dr Reddy's Laboratories Ltd"
$ws.Range("B3").Value = "Dr Reddy's Laboratories Ltd"

# Update numeric values on rows 2 and 3.
$ws.Range("A2").Value = 3213113113
$ws.Range("A3").Value = 1333333333
$ws.Range("C3").Value = 50

# Remove row 4 entirely (was A4=1234567845, B4="Dr Reddy's Laboratories Ltd.", C4=200).
$ws.Range("A4:C4").EntireRow.Delete()

# Match the selection left by the edit (B3 selected).
$ws.Range("B3").Select()
